# This script updates the cryptocurrency price table (columns B-E, rows 2-51)
# to reflect refreshed values from a scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.607.92"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +1.79%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.792.59"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.11"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.97%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.23"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +7.77%  "

# Row 9
$ws.Range("E9").Value = "  +0.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0691"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +3.88%  "

# Row 11
$ws.Range("E11").Value = "  +1.61%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.050.54"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.62%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.90"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +7.92%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.772.09"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.70%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.633.61"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +1.98%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.630"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +1.17%  "

# Row 17
$ws.Range("E17").Value = "  +2.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.47"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.56"
$ws.Range("D19").NumberFormat = "General"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0790"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +6.78%  "

# Row 21
$ws.Range("E21").Value = "  +0.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.43"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +1.03%  "

# Row 23
$ws.Range("E23").Value = "  -0.29%  "

# Row 24
$ws.Range("E24").Value = "  +0.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.35"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.62%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.38"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.70%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.09"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +1.72%  "

# Row 28
$ws.Range("E28").Value = "  -0.10%  "

# Row 29
$ws.Range("E29").Value = "  +0.16%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.75"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -1.80%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0515"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +0.25%  "

# Row 32
$ws.Range("E32").Value = "  +0.13%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.56"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.38%  "

# Row 34
$ws.Range("E34").Value = "  +2.34%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.436.19"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -3.97%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.43%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0191"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +2.91%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.635"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.77"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.56%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +3.17%  "

# Row 41
$ws.Range("E41").Value = "  +0.10%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.905"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +1.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.07"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +1.05%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.02"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +5.54%  "

# Row 45
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.06"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.39%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0494"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -3.52%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.948.67"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.83%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.00"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.35%  "

# Row 49
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.20%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.49"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +6.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.66"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -3.36%  "
